# Apply the grade-sheet update: fill in newly-graded "第10章作业" (col P)
# and "第5章作业" (col R) scores (and a couple of late entries in other
# columns) for the students whose homework/quiz results just came in,
# then leave the view scrolled/selected on the last-touched cell (R26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Score([int]$row, [int]$col, [double]$value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# --- Column I (第六章作业) late entry -------------------------------------
Set-Score 91 9 5

# --- Column O (第9章作业) late entries -------------------------------------
Set-Score 21 15 5
Set-Score 23 15 5
Set-Score 65 15 3

# --- Column P (第10章作业) new grades ---------------------------------------
Set-Score 2  16 5
Set-Score 3  16 5
Set-Score 7  16 3
Set-Score 9  16 5
Set-Score 11 16 5
Set-Score 14 16 3
Set-Score 16 16 5
Set-Score 19 16 5
Set-Score 21 16 5
Set-Score 22 16 5
Set-Score 23 16 5
Set-Score 24 16 5
Set-Score 26 16 3
Set-Score 29 16 5
Set-Score 30 16 5
Set-Score 43 16 5
Set-Score 45 16 5
Set-Score 50 16 5
Set-Score 51 16 5
Set-Score 53 16 5
Set-Score 55 16 5
Set-Score 57 16 5
Set-Score 62 16 5
Set-Score 65 16 3
Set-Score 66 16 5
Set-Score 67 16 5
Set-Score 70 16 5
Set-Score 75 16 5
Set-Score 80 16 5
Set-Score 91 16 5
Set-Score 92 16 5

# --- Column R (第5章作业) new grades -----------------------------------------
Set-Score 2  18 5
Set-Score 3  18 5
Set-Score 7  18 3
Set-Score 9  18 5
Set-Score 10 18 5
Set-Score 11 18 5
Set-Score 14 18 3
Set-Score 16 18 5
Set-Score 19 18 5
Set-Score 21 18 5
Set-Score 22 18 5
Set-Score 23 18 5
Set-Score 24 18 5
Set-Score 26 18 3
Set-Score 27 18 5
Set-Score 28 18 3
Set-Score 29 18 5
Set-Score 30 18 5
Set-Score 43 18 5
Set-Score 44 18 5
Set-Score 45 18 5
Set-Score 50 18 5
Set-Score 51 18 5
Set-Score 52 18 5
Set-Score 53 18 5
Set-Score 55 18 5
Set-Score 57 18 5
Set-Score 58 18 5
Set-Score 62 18 5
Set-Score 65 18 3
Set-Score 66 18 5
Set-Score 67 18 5
Set-Score 70 18 5
Set-Score 72 18 5
Set-Score 78 18 5
Set-Score 80 18 5
Set-Score 83 18 5
Set-Score 91 18 5
Set-Score 92 18 5

# --- Restore the view: scrolled so row 11 is at the top, with R26 selected --
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 1
[void]$ws.Range("R26").Select()
